# Applies the "GitHub Actions cryptos list" refresh described by the diff.
# Most edits are plain text overwrites (coin name / link / % change).
# A handful of "Price" cells are digit strings that Excel would otherwise
# auto-coerce to Number (dropping formatting like trailing zeros), so those
# are forced to Text via NumberFormat "@" and the format is cleared again
# right after so the cell style is left exactly as it started (General).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '58.933.21'
$ws.Range("E2").Value = '  +0.25%  '
# Row 3
$ws.Range("D3").Value = '2.502.54'
$ws.Range("E3").Value = '  +0.38%  '
# Row 4
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
Set-TextValue "D5" '532.37'
$ws.Range("E5").Value = '  -0.33%  '
# Row 6
Set-TextValue "D6" '135.28'
$ws.Range("E6").Value = '  -0.76%  '
# Row 7
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  +0.24%  '
# Row 8
$ws.Range("E8").Value = '  +0.49%  '
# Row 9
$ws.Range("E9").Value = '  +0.43%  '
# Row 10
$ws.Range("E10").Value = '  -1.10%  '
# Row 11
Set-TextValue "D11" '5.39'
$ws.Range("E11").Value = '  +0.64%  '
# Row 12
Set-TextValue "D12" '0.345'
$ws.Range("E12").Value = '  -0.48%  '
# Row 13
$ws.Range("D13").Value = '2.946.79'
$ws.Range("E13").Value = '  +0.21%  '
# Row 14
$ws.Range("D14").Value = '58.847.33'
$ws.Range("E14").Value = '  +0.16%  '
# Row 15
Set-TextValue "D15" '22.74'
$ws.Range("E15").Value = '  -1.62%  '
# Row 17
$ws.Range("D17").Value = '2.500.63'
$ws.Range("E17").Value = '  +0.05%  '
# Row 18
$ws.Range("E18").Value = '  -0.07%  '
# Row 19
Set-TextValue "D19" '4.23'
$ws.Range("E19").Value = '  -0.18%  '
# Row 20
Set-TextValue "D20" '323.58'
$ws.Range("E20").Value = '  -0.09%  '
# Row 21
$ws.Range("E21").Value = '  -0.06%  '
# Row 22
Set-TextValue "D22" '5.92'
$ws.Range("E22").Value = '  +1.17%  '
# Row 23
Set-TextValue "D23" '64.93'
$ws.Range("E23").Value = '  +0.64%  '
# Row 24
Set-TextValue "D24" '0.419'
$ws.Range("E24").Value = '  +0.29%  '
# Row 25
Set-TextValue "D25" '0.163'
$ws.Range("E25").Value = '  -0.72%  '
# Row 26
$ws.Range("E26").Value = '  +0.39%  '
# Row 27
Set-TextValue "D27" '7.50'
$ws.Range("E27").Value = '  -0.85%  '
# Row 28
$ws.Range("D28").Value = '0.0₃0761'
$ws.Range("E28").Value = '  -1.24%  '
# Row 29
Set-TextValue "D29" '6.45'
$ws.Range("E29").Value = '  -4.19%  '
# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D30" '1.74'
$ws.Range("E30").Value = '  -1.26%  '
# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D31" '168.93'
$ws.Range("E31").Value = '  +0.65%  '
# Row 32
$ws.Range("E32").Value = '  +0.11%  '
# Row 33
$ws.Range("E33").Value = '  -3.32%  '
# Row 34
Set-TextValue "D34" '18.32'
$ws.Range("E34").Value = '  -1.15%  '
# Row 35
$ws.Range("E35").Value = '  -4.01%  '
# Row 36
Set-TextValue "D36" '4.02'
$ws.Range("E36").Value = '  -1.25%  '
# Row 37
$ws.Range("E37").Value = '  -3.00%  '
# Row 38
$ws.Range("E38").Value = '  -1.14%  '
# Row 39
Set-TextValue "D39" '0.795'
$ws.Range("E39").Value = '  -3.17%  '
# Row 40
Set-TextValue "D40" '280.47'
$ws.Range("E40").Value = '  +0.66%  '
# Row 41
$ws.Range("E41").Value = '  +0.38%  '
# Row 42
$ws.Range("E42").Value = '  -0.24%  '
# Row 43
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D43" '10.91'
$ws.Range("E43").Value = '  +0.30%  '
# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D44" '4.97'
$ws.Range("E44").Value = '  -5.06%  '
# Row 45
Set-TextValue "D45" '129.10'
$ws.Range("E45").Value = '  +0.49%  '
# Row 46
$ws.Range("E46").Value = '  -0.85%  '
# Row 47
Set-TextValue "D47" '0.0499'
$ws.Range("E47").Value = '  -2.83%  '
# Row 48
$ws.Range("E48").Value = '  -1.19%  '
# Row 49
Set-TextValue "D49" '17.21'
$ws.Range("E49").Value = '  -0.55%  '
# Row 50
$ws.Range("D50").Value = '1.747.04'
$ws.Range("E50").Value = '  -1.10%  '
# Row 51
$ws.Range("E51").Value = '  -0.62%  '
